$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

# Row 2
$ws.Cells.Item(2, 4).Value = "30.244.59"
$ws.Cells.Item(2, 5).Value = "  -0.40%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.862.21"
$ws.Cells.Item(3, 5).Value = "  -0.52%  "

# Row 4
Set-TextValue $ws.Cells.Item(4, 4) "1.001"
$ws.Cells.Item(4, 5).Value = "  +0.01%  "

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) "236.76"
$ws.Cells.Item(5, 5).Value = "  +0.54%  "

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) "1.001"
$ws.Cells.Item(6, 5).Value = "  -0.03%  "

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) "0.4726"
$ws.Cells.Item(7, 5).Value = "  +1.18%  "

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) "0.2902"
$ws.Cells.Item(8, 5).Value = "  +2.07%  "

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) "0.06550"
$ws.Cells.Item(9, 5).Value = "  -0.25%  "

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) "21.83"
$ws.Cells.Item(10, 5).Value = "  +2.38%  "

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) "0.07933"
$ws.Cells.Item(11, 5).Value = "  +0.11%  "

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) "97.83"
$ws.Cells.Item(12, 5).Value = "  +0.34%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "1.865.76"
$ws.Cells.Item(13, 5).Value = "  -0.35%  "

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) "5.162"
$ws.Cells.Item(14, 5).Value = "  +0.53%  "

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) "0.6823"
$ws.Cells.Item(15, 5).Value = "  +0.85%  "

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) "267.30"
$ws.Cells.Item(16, 5).Value = "  -5.03%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "30.231.19"
$ws.Cells.Item(17, 5).Value = "  -0.42%  "

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) "13.74"
$ws.Cells.Item(18, 5).Value = "  +8.34%  "

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) "0.9999"
$ws.Cells.Item(19, 5).Value = "  -0.02%  "

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) "0.000007406"
$ws.Cells.Item(20, 5).Value = "  +1.41%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "2.107.97"
$ws.Cells.Item(21, 5).Value = "  +0.26%  "

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) "5.313"
$ws.Cells.Item(22, 5).Value = "  -3.82%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -0.01%  "

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) "6.188"
$ws.Cells.Item(24, 5).Value = "  -0.28%  "

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) "167.59"
$ws.Cells.Item(25, 5).Value = "  +1.57%  "

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) "9.230"
$ws.Cells.Item(26, 5).Value = "  -0.50%  "

# Row 27
Set-TextValue $ws.Cells.Item(27, 4) "18.94"
$ws.Cells.Item(27, 5).Value = "  -1.14%  "

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) "1.966"
$ws.Cells.Item(28, 5).Value = "  +1.18%  "

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) "1.394"
$ws.Cells.Item(29, 5).Value = "  +1.51%  "

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) "0.09874"
$ws.Cells.Item(30, 5).Value = "  +1.56%  "

# Row 31
Set-TextValue $ws.Cells.Item(31, 4) "4.376"
$ws.Cells.Item(31, 5).Value = "  -1.36%  "

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) "1.472"
$ws.Cells.Item(32, 5).Value = "  -0.41%  "

# Row 33
Set-TextValue $ws.Cells.Item(33, 4) "4.060"
$ws.Cells.Item(33, 5).Value = "  -1.36%  "

# Row 34
Set-TextValue $ws.Cells.Item(34, 4) "0.04716"
$ws.Cells.Item(34, 5).Value = "  +0.59%  "

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) "1.129"
$ws.Cells.Item(35, 5).Value = "  +1.10%  "

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) "0.7042"
$ws.Cells.Item(36, 5).Value = "  -0.27%  "

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) "2.707"
$ws.Cells.Item(37, 5).Value = "  -0.28%  "

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) "0.01880"
$ws.Cells.Item(38, 5).Value = "  +0.98%  "

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) "2.608"
$ws.Cells.Item(39, 5).Value = "  +2.46%  "

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) "6.265"
$ws.Cells.Item(40, 5).Value = "  -0.88%  "

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) "74.30"
$ws.Cells.Item(41, 5).Value = "  +1.33%  "

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) "1.942"
$ws.Cells.Item(42, 5).Value = "  -0.07%  "

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) "0.8454"
$ws.Cells.Item(43, 5).Value = "  -0.42%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "TheSandbox"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Cells.Item(44, 4) "0.4171"
$ws.Cells.Item(44, 5).Value = "  -0.47%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "PaxDollar"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Cells.Item(45, 4) "1.000"
$ws.Cells.Item(45, 5).Value = "  -0.11%  "

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) "103.52"
$ws.Cells.Item(46, 5).Value = "  -0.31%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "Aptos"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Cells.Item(47, 4) "7.181"
$ws.Cells.Item(47, 5).Value = "  -0.64%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "Maker"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Cells.Item(48, 4) "951.54"
$ws.Cells.Item(48, 5).Value = "  +2.06%  "

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) "9.214"
$ws.Cells.Item(49, 5).Value = "  +0.71%  "

# Row 50
Set-TextValue $ws.Cells.Item(50, 4) "34.16"
$ws.Cells.Item(50, 5).Value = "  -0.01%  "

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) "0.05655"
$ws.Cells.Item(51, 5).Value = "  +0.28%  "
